# Update the crypto price/volume table to reflect the latest scrape.
# Several "Price" column values look numeric (e.g. "603.87", "1.00") but must
# remain plain text exactly as scraped, so we apostrophe-prefix those assignments
# to force Excel to store them as text instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.196.49'
$ws.Range("D3").Value = '3.619.78'
$ws.Range("E3").Value = '  +3.57%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'" + '603.87'
$ws.Range("D6").Value = "'" + '196.09'
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("D7").Value = "'" + '0.627'
$ws.Range("E7").Value = '  +0.96%  '
$ws.Range("D8").Value = "'" + '1.00'
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").Value = "'" + '0.650'
$ws.Range("E10").Value = '  +0.15%  '
$ws.Range("D11").Value = "'" + '53.91'
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("E12").Value = '  +2.35%  '
$ws.Range("D13").Value = "'" + '9.55'
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("D14").Value = '4.202.60'
$ws.Range("E14").Value = '  +3.62%  '
$ws.Range("D15").Value = "'" + '13.19'
$ws.Range("E15").Value = '  +5.32%  '
$ws.Range("D16").Value = "'" + '591.99'
$ws.Range("E16").Value = '  -2.30%  '
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").Value = "'" + '19.26'
$ws.Range("E17").Value = '  +1.81%  '
$ws.Range("D18").Value = '70.428.83'
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.623.67'
$ws.Range("E19").Value = '  +3.69%  '
$ws.Range("E20").Value = '  +1.47%  '
$ws.Range("D21").Value = "'" + '0.997'
$ws.Range("E21").Value = '  +0.94%  '
$ws.Range("D22").Value = "'" + '17.68'
$ws.Range("E22").Value = '  -1.86%  '
$ws.Range("E23").Value = '  +1.28%  '
$ws.Range("D24").Value = "'" + '102.13'
$ws.Range("E24").Value = '  -2.12%  '
$ws.Range("D25").Value = "'" + '4.62'
$ws.Range("E25").Value = '  +1.44%  '
$ws.Range("E26").Value = '  -0.47%  '
$ws.Range("D27").Value = "'" + '10.81'
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("D28").Value = "'" + '9.62'
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("D29").Value = "'" + '34.06'
$ws.Range("E29").Value = '  +1.78%  '
$ws.Range("D30").Value = "'" + '4.79'
$ws.Range("E30").Value = '  +7.25%  '
$ws.Range("D31").Value = "'" + '7.17'
$ws.Range("E31").Value = '  +1.51%  '
$ws.Range("D32").Value = "'" + '12.35'
$ws.Range("E32").Value = '  -1.91%  '
$ws.Range("D33").Value = "'" + '0.117'
$ws.Range("E33").Value = '  +2.59%  '
$ws.Range("D34").Value = "'" + '0.0' + [char]0x2083 + '0906'
$ws.Range("E34").Value = '  +12.86%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("D36").Value = '3.918.41'
$ws.Range("E36").Value = '  +5.01%  '
$ws.Range("D38").Value = "'" + '529.04'
$ws.Range("E38").Value = '  +5.56%  '
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("D40").Value = "'" + '37.51'
$ws.Range("E40").Value = '  +2.74%  '
$ws.Range("D41").Value = "'" + '0.392'
$ws.Range("E41").Value = '  +0.77%  '
$ws.Range("D42").Value = "'" + '3.54'
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("E43").Value = '  -0.85%  '
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("D45").Value = "'" + '2.86'
$ws.Range("E45").Value = '  +1.89%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = "'" + '0.141'
$ws.Range("E46").Value = '  +1.04%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = "'" + '3.35'
$ws.Range("E47").Value = '  +0.75%  '
$ws.Range("D48").Value = "'" + '8.62'
$ws.Range("E48").Value = '  -0.96%  '
$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D49").Value = "'" + '0.000258'
$ws.Range("E49").Value = '  +7.74%  '
$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").Value = "'" + '1.00'
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").Value = "'" + '1.33'
$ws.Range("E51").Value = '  +5.06%  '
